$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for the team record (Wins/Losses/Ties), columns AD:AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold font, border, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every player row (rows 2-55): 82 wins, 80 losses, 0 ties
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 82   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 80   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
